$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.743.69'
$ws.Range('E2').Value = '  -0.97%  '
$ws.Range('D3').Value = '2.370.93'
$ws.Range('E3').Value = '  +1.23%  '
$ws.Range('E4').Value = '  -0.40%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '332.50'
$ws.Range('E5').Value = '  +6.29%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '100.82'
$ws.Range('E6').Value = '  -7.83%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.639'
$ws.Range('E7').Value = '  -0.22%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.631'
$ws.Range('E9').Value = '  +0.22%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.01'
$ws.Range('E10').Value = '  -6.38%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0924'
$ws.Range('E11').Value = '  -1.30%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.49'
$ws.Range('E12').Value = '  -4.29%  '
$ws.Range('E13').Value = '  -3.06%  '
$ws.Range('E14').Value = '  +0.29%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '16.46'
$ws.Range('E15').Value = '  +1.27%  '
$ws.Range('D16').Value = '2.724.75'
$ws.Range('E16').Value = '  +0.72%  '
$ws.Range('D17').Value = '2.365.56'
$ws.Range('E17').Value = '  -4.80%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.96'
$ws.Range('E18').Value = '  +9.64%  '
$ws.Range('D19').Value = '42.663.05'
$ws.Range('E19').Value = '  -1.15%  '
$ws.Range('E20').Value = '  -1.26%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '3.79'
$ws.Range('E21').Value = '  +10.26%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '76.21'
$ws.Range('E22').Value = '  +1.28%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '269.78'
$ws.Range('E23').Value = '  +6.74%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.32'
$ws.Range('E24').Value = '  -10.38%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.08'
$ws.Range('E25').Value = '  +10.92%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.52'
$ws.Range('E27').Value = '  -4.05%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '23.17'
$ws.Range('E28').Value = '  +3.67%  '
$ws.Range('E29').Value = '  -2.58%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '176.59'
$ws.Range('E30').Value = '  +1.17%  '
$ws.Range('E31').Value = '  -2.84%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0903'
$ws.Range('E32').Value = '  -3.12%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '35.44'
$ws.Range('E33').Value = '  -9.87%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.18'
$ws.Range('E34').Value = '  +3.22%  '
$ws.Range('E35').Value = '  +0.52%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.62'
$ws.Range('E36').Value = '  -6.75%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.99'
$ws.Range('E37').Value = '  +11.54%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0360'
$ws.Range('E38').Value = '  -4.51%  '
$ws.Range('B39').Value = 'NEARProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.84'
$ws.Range('E39').Value = '  -7.18%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.105'
$ws.Range('E40').Value = '  +1.41%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.54'
$ws.Range('E41').Value = '  +4.60%  '
$ws.Range('E42').Value = '  +1.45%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '70.35'
$ws.Range('E43').Value = '  -3.18%  '
$ws.Range('E44').Value = '  -0.21%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '118.61'
$ws.Range('E45').Value = '  +7.32%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '91.42'
$ws.Range('E46').Value = '  +30.84%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '11.93'
$ws.Range('E47').Value = '  -7.10%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.49'
$ws.Range('E48').Value = '  -2.32%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.18'
$ws.Range('E49').Value = '  -1.34%  '
$ws.Range('E50').Value = '  -2.28%  '
$ws.Range('D51').Value = '1.571.57'
$ws.Range('E51').Value = '  +5.53%  '
